# AdminSection Module datepickerUtil added
# Update the ECM test-data credentials used by the Admin Section module
# from the "Ecm04" test user to the "Ecm06" test user, and move the
# active selection to E32 as left by the author after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Update ECM_FirstName
$ws.Range("B32").Value = "Ecm06"

# ECM_LastName (B33) stays "Testuser" - no change in value

# Update ECM_UserName
$ws.Range("B34").Value = "Ecm06_testuser"

# Update ECM_EmailId
$ws.Range("B35").Value = "testuser006@test.com"

# Update EditUserName value to match the new ECM username
$ws.Range("B41").Value = "Ecm06_testuser"

# Move the active selection/cell as seen in the saved workbook
$ws.Range("E32").Select()
